# Standards - Automation - 12/24/2025
#
# 1) Bump the "Last id recorded" counter from 2 to 3.
# 2) Record a new procedure entry ("3: How to cleaning screen") as a new
#    numbered list item right after the existing "How to make Test Case"
#    entry, matching the existing list's paragraph style / numbering.

$d = $word.ActiveDocument

# --- 1) Last id recorded: 2 -> 3 -------------------------------------------
$d.Content.Find.Execute("Last id recorded: 2", $true, $false, $false, $false,
                         $false, $true, 1, $false, "Last id recorded: 3", 2) | Out-Null

# --- 2) Add "3: How to cleaning screen" after "2: How to make Test Case" ---
# Locate the existing "How to make Test Case" list item by its text so the
# edit is resilient to the exact paragraph index.
$targetIndex = 0
$i = 0
foreach ($para in $d.Paragraphs) {
    $i = $i + 1
    if ($para.Range.Text -like "*How to make Test Case*") {
        $targetIndex = $i
    }
}

$target = $d.Paragraphs($targetIndex)

# Inserting a paragraph break after the existing item creates a new
# paragraph that inherits its style (ListParagraph) and numbering
# (ilvl 0 / numId 1), matching the rest of the list.
$target.Range.InsertParagraphAfter()

$newItem = $d.Paragraphs($targetIndex + 1)
$newItem.Range.Text = "3: "

$endOfNewItem = $d.Paragraphs($targetIndex + 1).Range
$endOfNewItem.Collapse(0)
$endOfNewItem.InsertAfter("How to cleaning screen")
